$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1387.3478
$ws.Range("I113").Value = 1439.1428
$ws.Range("J113").Value = 1306.7778
$ws.Range("K113").Value = 1439.1428
$ws.Range("L113").Value = 1306.7778
$ws.Range("M113").Value = 1814.8572
$ws.Range("N113").Value = -7814.7778
$ws.Range("H132").Value = 31761.727
$ws.Range("I132").Value = 38353.43
$ws.Range("K132").Value = 115060.29
$ws.Range("M132").Value = -112530.29
$ws.Range("H139").Value = 89697.664
$ws.Range("J139").Value = 89697.664
$ws.Range("L139").Value = 89697.664
$ws.Range("N139").Value = -99977.664
$ws.Range("H140").Value = 148999.33
$ws.Range("J140").Value = 148999.33
$ws.Range("L140").Value = 148999.33
$ws.Range("N140").Value = -159359.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1112.1621
$ws.Range("I2").Value = 1078.5483
$ws.Range("K2").Value = 1078.5483
$ws.Range("M2").Value = -965.5482999999999
$ws.Range("H63").Value = 10500.625
$ws.Range("I63").Value = 5005
$ws.Range("J63").Value = 11285.714
$ws.Range("K63").Value = 5005
$ws.Range("L63").Value = 11285.714
$ws.Range("M63").Value = -4319
$ws.Range("N63").Value = -12657.714
$ws.Range("H66").Value = 10500.625
$ws.Range("I66").Value = 5005
$ws.Range("J66").Value = 11285.714
$ws.Range("K66").Value = 25025
$ws.Range("L66").Value = 56428.57
$ws.Range("M66").Value = -21593
$ws.Range("N66").Value = -63292.57
$ws.Range("H101").Value = 76754.71000000001
$ws.Range("J101").Value = 76754.71000000001
$ws.Range("L101").Value = 76754.71000000001
$ws.Range("N101").Value = -83244.71000000001
$ws.Range("H110").Value = 1081.0555
$ws.Range("I110").Value = 1081.0555
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1081.0555
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = 963.9445000000001
$ws.Range("H116").Value = 1112.1621
$ws.Range("I116").Value = 1078.5483
$ws.Range("K116").Value = 1078.5483
$ws.Range("M116").Value = 1215.4517

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1112.1621
$ws.Range("I3").Value = 1078.5483
$ws.Range("K3").Value = 1078.5483
$ws.Range("M3").Value = -964.5482999999999
$ws.Range("H80").Value = 4331.6665
$ws.Range("I80").Value = 8690.166999999999
$ws.Range("J80").Value = 2152.4167
$ws.Range("K80").Value = 8690.166999999999
$ws.Range("L80").Value = 2152.4167
$ws.Range("M80").Value = -7692.166999999999
$ws.Range("N80").Value = -4148.4167
$ws.Range("H83").Value = 4331.6665
$ws.Range("I83").Value = 8690.166999999999
$ws.Range("J83").Value = 2152.4167
$ws.Range("K83").Value = 43450.835
$ws.Range("L83").Value = 10762.0835
$ws.Range("M83").Value = -38458.835
$ws.Range("N83").Value = -20746.0835
$ws.Range("H86").Value = 3318
$ws.Range("I86").Value = 2761.25
$ws.Range("J86").Value = 3874.75
$ws.Range("K86").Value = 2761.25
$ws.Range("L86").Value = 3874.75
$ws.Range("M86").Value = -1638.25
$ws.Range("N86").Value = -6120.75
$ws.Range("H89").Value = 3318
$ws.Range("I89").Value = 2761.25
$ws.Range("J89").Value = 3874.75
$ws.Range("K89").Value = 13806.25
$ws.Range("L89").Value = 19373.75
$ws.Range("M89").Value = -8190.25
$ws.Range("N89").Value = -30605.75
$ws.Range("H94").Value = 1529.625
$ws.Range("I94").Value = 1147.4
$ws.Range("K94").Value = 1147.4
$ws.Range("M94").Value = -696.4000000000001
$ws.Range("H99").Value = 1679.4
$ws.Range("I99").Value = 1299.4286
$ws.Range("K99").Value = 1299.4286
$ws.Range("M99").Value = 198.5714
$ws.Range("H105").Value = 2032.9166
$ws.Range("J105").Value = 2815.5
$ws.Range("L105").Value = 2815.5
$ws.Range("N105").Value = -6309.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5653.475
$ws.Range("I31").Value = 1574.3158
$ws.Range("J31").Value = 9344.143
$ws.Range("K31").Value = 1574.3158
$ws.Range("L31").Value = 9344.143
$ws.Range("M31").Value = -1279.3158
$ws.Range("N31").Value = -9934.143
$ws.Range("H34").Value = 5653.475
$ws.Range("I34").Value = 1574.3158
$ws.Range("J34").Value = 9344.143
$ws.Range("K34").Value = 1574.3158
$ws.Range("L34").Value = 9344.143
$ws.Range("M34").Value = -1372.3158
$ws.Range("N34").Value = -9748.143
$ws.Range("H99").Value = 7025
$ws.Range("I99").Value = 6321.4287
$ws.Range("J99").Value = 7728.5713
$ws.Range("K99").Value = 6321.4287
$ws.Range("L99").Value = 7728.5713
$ws.Range("M99").Value = -4823.4287
$ws.Range("N99").Value = -10724.5713
$ws.Range("H126").Value = 7025
$ws.Range("I126").Value = 6321.4287
$ws.Range("J126").Value = 7728.5713
$ws.Range("K126").Value = 18964.2861
$ws.Range("L126").Value = 23185.7139
$ws.Range("M126").Value = -16494.2861
$ws.Range("N126").Value = -28125.7139
$ws.Range("H132").Value = 9616939
$ws.Range("I132").Value = 1630.9131
$ws.Range("K132").Value = 4892.7393
$ws.Range("M132").Value = -2362.7393
$ws.Range("H134").Value = 2139.4783
$ws.Range("I134").Value = 2060.4
$ws.Range("K134").Value = 6181.200000000001
$ws.Range("M134").Value = -3646.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2183.923
$ws.Range("J97").Value = 4361.0835
$ws.Range("L97").Value = 4361.0835
$ws.Range("N97").Value = -5353.0835

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H7").Value = 7176.2856
$ws.Range("I7").Value = 6705.6665
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 6705.6665
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -6593.6665
$ws.Range("N7").Value = -10224
$ws.Range("H46").Value = 2879.6296
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 2910
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 2910
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -3286
$ws.Range("H100").Value = 6756.7393
$ws.Range("I100").Value = 2422.0557
$ws.Range("J100").Value = 22361.6
$ws.Range("K100").Value = 2422.0557
$ws.Range("L100").Value = 22361.6
$ws.Range("M100").Value = -1881.0557
$ws.Range("N100").Value = -23443.6
$ws.Range("H126").Value = 7176.2856
$ws.Range("I126").Value = 6705.6665
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 20116.9995
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -17646.9995
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 911918.75
$ws.Range("I132").Value = 1194062.1
$ws.Range("K132").Value = 3582186.3
$ws.Range("M132").Value = -3579656.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1133.375
$ws.Range("I100").Value = 836.25
$ws.Range("K100").Value = 1672.5
$ws.Range("M100").Value = -1131.5
$ws.Range("H122").Value = 3145.2
$ws.Range("I122").Value = 2331.3572
$ws.Range("K122").Value = 6994.071599999999
$ws.Range("M122").Value = -4544.071599999999
